# Convert the four Word "begin/instrText/end" field codes used as M2Doc
# tags ( m: ... ) into plain literal text runs written with the
# { m: ... } curly-brace syntax expected by the
# TokenIteratorFieldRewriterSplit parser.
#
# Each field is deleted (which removes its begin/instrText/end runs) and
# replaced by a plain w:t run holding the equivalent "{...}" text, while
# keeping the paragraph's run formatting (w:lang) intact.
#
# Fields are processed from the last one in the document back to the
# first so that earlier field/paragraph indices and Range offsets stay
# valid while later ones are being rewritten.

$d = $word.ActiveDocument

# --- Paragraph 6: " m:endtemplate " -> "{m:endtemplate}" ------------------
$p6 = $d.Paragraphs.Item(6)
$f6 = $d.Fields.Item($d.Fields.Count)
$f6.Delete()
$p6.Range.Text = "{m:endtemplate}"
$p6.Range.LanguageID = "en-US"

# --- Paragraph 5: " m: a + a " -> "{m: a + a}" -----------------------------
$p5 = $d.Paragraphs.Item(5)
$f5 = $d.Fields.Item($d.Fields.Count)
$f5.Delete()
$p5.Range.Text = "{m: a + a}"
$p5.Range.LanguageID = "en-US"

# --- Paragraph 4: " m:templa" + bookmark + " myTemplate(a:Integer) " ------
# -> "{m:templa" + bookmark + " myTemplate(a:Integer)}"
# The _GoBack bookmark that sits between the two instrText runs has to be
# preserved at the same logical split point (right after "templa").
$p4 = $d.Paragraphs.Item(4)
$p4Start = $p4.Range.Start
$f4 = $d.Fields.Item($d.Fields.Count)
$f4.Delete()
$p4.Range.Text = "{m:templa myTemplate(a:Integer)}"
$p4.Range.LanguageID = "en-US"
$bmPos = $p4Start + 9
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Paragraph 2: " m: 2.myTemplate() " -> "{m: 2.myTemplate()}" ----------
$p2 = $d.Paragraphs.Item(2)
$f2 = $d.Fields.Item($d.Fields.Count)
$f2.Delete()
$p2.Range.Text = "{m: 2.myTemplate()}"
$p2.Range.LanguageID = "en-US"
